# Applies:
#  1. Merge the "Project Proposal and Outline" run and the following
#     single-space run into one run reading "Project Proposal and Outline ".
#  2. Strike through the "Sharpe ratio" and "Beta" list paragraphs
#     (run text + paragraph mark).

$d = $word.ActiveDocument

# --- 1. Merge "Project Proposal and Outline" + " " into a single run ---
$d.Content.Find.Execute(
    "Project Proposal and Outline ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Project Proposal and Outline ", 2) | Out-Null

# --- 2. Strike through "Sharpe ratio" and "Beta" paragraphs ---
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.Trim()
    if ($txt -eq "Sharpe ratio" -or $txt -eq "Beta") {
        $p.Range.Font.StrikeThrough = 1
    }
}
